# Add two new columns I (I0) and J (IF) to the sheet, mirroring the style
# of the existing header cells and filling in the per-row numeric data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style used by the other header cells (e.g. H1) onto the
# two new header cells so they match (bold, bordered, centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data values (rows 2-41) ---
$iValues = @(8,9,6,8,7,6,7,8,7,5,6,7,5,7,4,6,5,6,7,6,5,3,9,2,6,7,5,6,5,8,5,6,5,6,9,5,8,8,3,9)
$jValues = @(8,9,6,8,7,6,7,8,8,5,6,7,5,8,4,6,5,6,8,6,6,4,9,2,6,7,5,6,6,8,6,7,6,6,9,5,8,8,3,9)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
